# Applies the diff: splits several runs so that embedded pseudo-XML
# markup (<exp>...</exp>, <ill/>) gets its own "tag" run formatting
# (Courier New / grey / small), renames one <m>...</m> pair to
# <pa>...</pa>, and gives the <ill/> tag its own Courier/blue styling.

$d = $word.ActiveDocument

# OLE/COM Font.Color is 0x00BBGGRR (reverse byte order from the usual
# HTML/OOXML 0xRRGGBB), so convert explicitly instead of hard-coding.
function RgbToOle($rgbHex) {
    $r = [Convert]::ToInt32($rgbHex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($rgbHex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($rgbHex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$colorGrey = RgbToOle "a9a9a9"
$colorBlue = RgbToOle "0000ff"

# Locate a (unique) run of text anywhere in the document body and
# return the Range covering exactly that text.
function Find-UniqueRange($needle) {
    $r = $d.Content
    $found = $r.Find.Execute($needle, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        throw "text not found: $needle"
    }
    return $d.Range($r.Start, $r.End)
}

# Style a sub-range as an inline "<tag>" marker: small grey Courier New.
function Style-ExpTag($range) {
    $range.Font.Name = "Courier New"
    $range.Font.Size = 7
    $range.Font.Color = $colorGrey
}

# Split "<before><exp><mid></exp><after>" (currently one run) into five
# runs: before / <exp> / mid / </exp> / after -- only the two tag runs
# get the new grey Courier New styling, the rest keep their inherited
# (untouched) formatting.
function Split-ExpRun($fullText, $before, $mid, $after) {
    $whole = Find-UniqueRange $fullText
    $start = $whole.Start

    $beforeEnd = $start + $before.Length
    $openStart = $beforeEnd
    $openEnd = $openStart + 5            # "<exp>"
    $midStart = $openEnd
    $midEnd = $midStart + $mid.Length
    $closeStart = $midEnd
    $closeEnd = $closeStart + 6          # "</exp>"
    $afterEnd = $closeEnd + $after.Length

    $openRange = $d.Range($openStart, $openEnd)
    if ($openRange.Text -ne "<exp>") { throw "unexpected open text: $($openRange.Text)" }
    Style-ExpTag $openRange

    $closeRange = $d.Range($closeStart, $closeEnd)
    if ($closeRange.Text -ne "</exp>") { throw "unexpected close text: $($closeRange.Text)" }
    Style-ExpTag $closeRange
}

Split-ExpRun "curieusem<exp>ent</exp> ceste pouldre Car un ℥ dicelle gectee sur" `
             "curieusem" "ent" " ceste pouldre Car un ℥ dicelle gectee sur"

Split-ExpRun "subtillem<exp>ent</exp> battu en lamines ou aultres ouvrages &" `
             "subtillem" "ent" " battu en lamines ou aultres ouvrages &"

Split-ExpRun "Le terme de loeuvre est de noeuf moys despuys la s<exp>ainc</exp>t Jehan" `
             "Le terme de loeuvre est de noeuf moys despuys la s" "ainc" "t Jehan"

Split-ExpRun "que la grene soict parfaictem<exp>ent</exp> meure le pied est si deseche du soleil" `
             "que la grene soict parfaictem" "ent" " meure le pied est si deseche du soleil"

# <m>Lin</m> -> <pa>Lin</pa> (only this one occurrence; other <m>...</m>
# pairs elsewhere in the document stay untouched).
$mOpen = Find-UniqueRange "<m>Lin</m>"
$mOpenStart = $mOpen.Start
$openTag = $d.Range($mOpenStart, $mOpenStart + 3)
if ($openTag.Text -ne "<m>") { throw "unexpected: $($openTag.Text)" }
$openTag.Text = "<pa>"

$closeTag = $d.Range($mOpenStart + 7, $mOpenStart + 11)
if ($closeTag.Text -ne "</m>") { throw "unexpected: $($closeTag.Text)" }
$closeTag.Text = "</pa>"

# f<ill/> -> split into "f" (unchanged) + "<ill/>" (new Courier/blue tag)
$fill = Find-UniqueRange "f<ill/>"
$fillStart = $fill.Start
$illRange = $d.Range($fillStart + 1, $fillStart + 7)
if ($illRange.Text -ne "<ill/>") { throw "unexpected: $($illRange.Text)" }
$illRange.Font.Name = "Courier New"
$illRange.Font.Size = 9
$illRange.Font.Color = $colorBlue
$illRange.Font.StrikeThrough = $false

Write-Output "done"
